$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.459.52"
$ws.Range("E2").Value = "  +4.60%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.816.24"
$ws.Range("E3").Value = "  +5.54%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "343.43"
$ws.Range("E5").Value = "  +2.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9990"
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3817"
$ws.Range("E7").Value = "  +3.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3503"
$ws.Range("E8").Value = "  +4.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "49.04"
$ws.Range("E9").Value = "  -0.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.234"
$ws.Range("E10").Value = "  +3.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07749"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.31"
$ws.Range("E13").Value = "  +10.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.606"
$ws.Range("E14").Value = "  +4.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.813.31"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.232"
$ws.Range("E17").Value = "  +3.70%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06719"
$ws.Range("E18").Value = "  +1.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "86.19"
$ws.Range("E19").Value = "  +4.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9989"
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.61"
$ws.Range("E21").Value = "  +7.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.574"
$ws.Range("E22").Value = "  +7.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.21"
$ws.Range("E23").Value = "  +1.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "27.453.02"
$ws.Range("E24").Value = "  +4.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.468"
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.670"
$ws.Range("E26").Value = "  +7.75%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.04"
$ws.Range("E27").Value = "  +14.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.471"
$ws.Range("E28").Value = "  +7.27%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "154.08"
$ws.Range("E29").Value = "  +1.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.019.69"
$ws.Range("E30").Value = "  +5.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "136.02"
$ws.Range("E31").Value = "  +4.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.345"
$ws.Range("E32").Value = "  +5.69%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.053"
$ws.Range("E33").Value = "  -1.21%  "
$ws.Range("E34").Value = "  +7.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08797"
$ws.Range("E35").Value = "  +3.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.692"
$ws.Range("E36").Value = "  -0.75%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.614"
$ws.Range("E37").Value = "  +4.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.6974"
$ws.Range("E38").Value = "  +12.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2265"
$ws.Range("E39").Value = "  +5.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02401"
$ws.Range("E40").Value = "  +2.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.06478"
$ws.Range("E41").Value = "  +3.41%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.961"
$ws.Range("E42").Value = "  +4.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.296"
$ws.Range("E43").Value = "  +5.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.83"
$ws.Range("E44").Value = "  +2.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6539"
$ws.Range("E45").Value = "  +10.18%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9986"
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.005"
$ws.Range("E47").Value = "  +3.87%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.176"
$ws.Range("E48").Value = "  +7.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "132.87"
$ws.Range("E49").Value = "  +3.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07325"
$ws.Range("E50").Value = "  +0.58%  "

# Row 51: Stacks -> Aave
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.48"
$ws.Range("E51").Value = "  +4.21%  "
